# Add a new sheet "ओन्ली सचिब" ("Secretary Only") as a simplified duplicate of
# the "ESTIMATE (3)" sheet: the computed-cost columns (start/end meter, elapsed
# time and the derived totals) are cleared, and the three signature-line rows
# near the bottom lose the "ward representative" column, sliding the
# remaining labels over (and dropping the old "ward technician" label).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the sheet, the copy is inserted right after the original and
# becomes the active sheet/tab automatically.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "ओन्ली सचिब"

# Clear the meter-reading / rate figures and every total derived from them.
$ws2.Range("C12:E12").ClearContents()
$ws2.Range("G12").ClearContents()
$ws2.Range("G13").ClearContents()
$ws2.Range("G15:G17").ClearContents()

# Row 20 ("...." placeholder signature lines): drop the C column, shifting
# its neighbour right, and drop the E column entirely.
$ws2.Range("C20").Copy($ws2.Range("D20"))
$ws2.Range("C20").Clear()
$ws2.Range("E20").Clear()

# Row 21 (blank underline cells beneath row 20) follows the same shift.
$ws2.Range("C21").Copy($ws2.Range("D21"))
$ws2.Range("C21").Clear()
$ws2.Range("E21").Clear()

# Row 22 (signature captions): ward-representative moves C->D, the blank
# spacer moves D->F, ward-secretary moves E->G (replacing the old
# ward-technician caption), then the now-vacated C/E cells are cleared.
$ws2.Range("E22").Copy($ws2.Range("G22"))
$ws2.Range("D22").Copy($ws2.Range("F22"))
$ws2.Range("C22").Copy($ws2.Range("D22"))
$ws2.Range("C22").Clear()
$ws2.Range("E22").Clear()

# Give the new sheet its own print area / print title rows (mirrors the
# ones already defined for "ESTIMATE (3)").
$ws2.PageSetup.PrintArea = "`$A`$1:`$H`$22"
$ws2.PageSetup.PrintTitleRows = "`$1:`$11"

# Restore the view the new sheet was left on.
$ws2.Activate()
$ws2.Range("D25").Select()
